$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 38
$ws1.Range("F4").Value = 133
$ws1.Range("F5").Value = 980
$ws1.Range("F6").Value = 395
$ws1.Range("F7").Value = 7602
$ws1.Range("F8").Value = 118
$ws1.Range("F9").Value = 174
$ws1.Range("F10").Value = 6698
$ws1.Range("F12").Value = 287
$ws1.Range("F13").Value = 4724
$ws1.Range("F17").Value = 4980
$ws1.Range("F18").Value = 1050
$ws1.Range("F19").Value = 283
$ws1.Range("F20").Value = 288
$ws1.Range("F21").Value = 389
$ws1.Range("F24").Value = 131
$ws1.Range("F25").Value = 84
$ws1.Range("F26").Value = 8605
$ws1.Range("F27").Value = 63
$ws1.Range("F29").Value = 33
$ws1.Range("F30").Value = 757
$ws1.Range("F33").Value = 65
$ws1.Range("F35").Value = 993
$ws1.Range("F37").Value = 1768
$ws1.Range("F39").Value = 1059
$ws1.Range("F41").Value = 4506
$ws1.Range("F44").Value = 128
$ws1.Range("F46").Value = 9
$ws1.Range("F47").Value = 874
$ws1.Range("F48").Value = 1182

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 32
$ws2.Range("F3").Value = 28
$ws2.Range("F17").Value = 877

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 38
$ws4.Range("F5").Value = 133
$ws4.Range("F6").Value = 28
$ws4.Range("F7").Value = 980
$ws4.Range("F8").Value = 395
$ws4.Range("F9").Value = 7602
$ws4.Range("F10").Value = 118
$ws4.Range("F11").Value = 174
$ws4.Range("F12").Value = 6698
$ws4.Range("F14").Value = 287
$ws4.Range("F15").Value = 4724
$ws4.Range("F19").Value = 4980
$ws4.Range("F20").Value = 1050
$ws4.Range("F21").Value = 283
$ws4.Range("F22").Value = 288
$ws4.Range("F23").Value = 389
$ws4.Range("F26").Value = 131
$ws4.Range("F27").Value = 84
$ws4.Range("F29").Value = 8605
$ws4.Range("F30").Value = 63
$ws4.Range("F32").Value = 33
$ws4.Range("F33").Value = 757
$ws4.Range("F35").Value = 65
$ws4.Range("F37").Value = 993
$ws4.Range("F38").Value = 1768
$ws4.Range("F40").Value = 1059
$ws4.Range("F42").Value = 4506
$ws4.Range("F45").Value = 128
$ws4.Range("F47").Value = 874
$ws4.Range("F48").Value = 1182
